$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 313, shifting the existing data rows (old 313:343) down to 314:344
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new weekly price record
$ws.Cells.Item(313, 1).Value = 7
$ws.Cells.Item(313, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(313, 3).Value = "Ñuble"
$ws.Cells.Item(313, 4).Value = 44858
$ws.Cells.Item(313, 5).Value = 16
$ws.Cells.Item(313, 6).Value = 100114013
$ws.Cells.Item(313, 7).Value = "Zanahoria"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 120
$ws.Cells.Item(313, 11).Value = 17000
$ws.Cells.Item(313, 12).Value = 18000
$ws.Cells.Item(313, 13).Value = 17500
$ws.Cells.Item(313, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(313, 15).Value = "Región de Ñuble"
$ws.Cells.Item(313, 16).Value = 875
$ws.Cells.Item(313, 17).Value = 20
$ws.Cells.Item(313, 18).Value = "Hortaliza"
